$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column H header + "ALL" marker on row 2 (mirrors existing 寄送對象/yes columns)
$ws.Range("H1").Value = "備註"
$ws.Range("H2").Value = "ALL"

# New row 3: another recipient entry, same shape as row 2
$ws.Range("A3").Value = "victor.hou@kingza.com.tw"
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:victor.hou@kingza.com.tw") | Out-Null
$ws.Range("A3").Style = $ws.Range("A2").Style
$ws.Range("B3").Value = "yes"
$ws.Range("G3").Value = "yes"
$ws.Range("H3").Value = "洪國瑋"
$ws.Range("H3").NumberFormat = "#,##0"
$ws.Range("H3").HorizontalAlignment = -4152

# Leave the selection where it ended up after filling the row
$ws.Range("G3").Select() | Out-Null
